$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 must stay text ("003"), not be auto-converted to the number 3
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "003"

$ws.Range("N2").Value = "2020-03-31 00:00:00"
$ws.Range("O2").Value = 15102255213.37
$ws.Range("P2").Value = 2995722165.58
$ws.Range("Q2").Value = 4488677584.99
$ws.Range("R2").Value = 26.5610242522
$ws.Range("S2").Value = 198432197.64
$ws.Range("T2").Value = 46.8735448525
$ws.Range("U2").Value = 2590930476.47
$ws.Range("V2").Value = 11.7780141547
$ws.Range("W2").Value = 8488037523.46
$ws.Range("X2").Value = 2343884081.2

# Y2 becomes blank/empty in the new data
$ws.Range("Y2").ClearContents()

$ws.Range("Z2").Value = 1644252450.75
$ws.Range("AA2").Value = 15.7331756688
$ws.Range("AB2").Value = 6614217689.91
$ws.Range("AC2").Value = 6.6402688928
$ws.Range("AD2").Value = 9.7531766726
$ws.Range("AE2").Value = 12.3077955664
$ws.Range("AF2").Value = 112.7119307023
$ws.Range("AG2").Value = 56.2037748902
